# Twitter "API Key" line: move the space that currently sits right after the
# manual line break (<w:br/>) so that it instead sits right after the "="
# sign, immediately before the key value.
#
#   before:  <break>" API Key =155C4FzRR7wjX6RtHlChscRgf"
#   after:   <break>"API Key = 155C4FzRR7wjX6RtHlChscRgf"

$d = $word.ActiveDocument

# Find the paragraph by its (unique) key value instead of a hard-coded
# paragraph index, so the script is resilient to unrelated structural
# differences elsewhere in the document.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*155C4FzRR7wjX6RtHlChscRgf*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $text = $r.Text

    # The paragraph starts with a manual line break (chr 11), immediately
    # followed by a single space - remove that space.
    $breakIdx = $text.IndexOf([char]11)
    if ($breakIdx -ge 0) {
        $spacePos = $r.Start + $breakIdx + 1
        $spaceRange = $d.Range($spacePos, $spacePos + 1)
        if ($spaceRange.Text -eq " ") {
            $spaceRange.Text = ""
        }
    }

    # Re-read the paragraph range/text since the delete above shifted
    # offsets, then insert a single space right after "=" (before the
    # key value), unless it is already there.
    $r2 = $target.Range
    $text2 = $r2.Text
    $eqIdx = $text2.IndexOf("=")
    if ($eqIdx -ge 0) {
        $afterEqPos = $r2.Start + $eqIdx + 1
        $afterEq = $d.Range($afterEqPos, $afterEqPos + 1)
        if ($afterEq.Text -ne " ") {
            $insertPoint = $d.Range($afterEqPos, $afterEqPos)
            $insertPoint.InsertAfter(" ")
        }
    }
}
